$d = $word.ActiveDocument

# Locate the last paragraph in the body (the "Si on doit faire une
# architecture MVC..." question) and position an insertion point right
# after it.
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n)
$r = $lastPara.Range
$r.Collapse(0)

# Three blank paragraphs.
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)

# Paragraph with the new bullet text, written as two separate runs
# ("- " then the rest), matching how it was originally typed.
$r.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$rp = $newPara.Range
$rp.Collapse(0)
$rp.InsertAfter("- Setup.py ou Frozenpython pour installer automatiquement les modules pythons puis exécuter le programme (création d’un .exe)")

# Force the "- " prefix to live in its own run distinct from the rest
# of the sentence (identical formatting, so toggle bold on/off on just
# that substring to split the run without leaving any visible change).
$full = $newPara.Range
$prefix = $d.Range($full.Start, $full.Start + 2)
$prefix.Font.Bold = 1
$prefix.Font.Bold = 0

# Trailing blank paragraph.
$rp2 = $newPara.Range
$rp2.Collapse(0)
$rp2.InsertParagraphAfter()

Write-Output "done"
